$wb = $excel.ActiveWorkbook

# --- Update the header/template row (row 3, columns B:I) on both sheets
# with data-binding placeholder text. ---
foreach ($ws in @($wb.Worksheets.Item(1), $wb.Worksheets.Item(2))) {
    $ws.Range("B3").Value = "&=dataSource.applicantName"
    $ws.Range("C3").Value = "&=dataSource.appType"
    $ws.Range("D3").Value = "&=dataSource.prePostAtr"
    $ws.Range("E3").Value = "&=dataSource.appStartDate"
    $ws.Range("F3").Value = "&=dataSource.appContent"
    $ws.Range("G3").Value = "&=dataSource.inputDate"
    $ws.Range("H3").Value = "&=dataSource.reflectionStatus"

    # I3 already carries a "quote prefix" cell style (it was formatted for
    # text that looks like a formula). Prefix the literal value with a
    # leading apostrophe, the same way a user types it in Excel, so the
    # quote-prefix formatting is preserved instead of being dropped.
    $ws.Range("I3").Value = "'&=dataSource.approvalStatusInquiry"
}

# --- Remove the print area defined for the first sheet (承認一覧). ---
$wb.Worksheets.Item(1).PageSetup.PrintArea = ""
